{"js": "// Applies the textual edits described by the commit diff:\n//  1. \"navigate the menus game\" -> \"navigate the game menus,\" (reworded + comma)\n//  2. \"menu selector\" -> \"level selector\"\n//  3. \"The Options menu\" -> \"The options menu\" (capitalization)\n//  4. add \", the menu will also be accessible through the pause menu\" before the final period\n//  5. \"will feel like\" -> \"could feel like\"\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replaceText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for '\" + searchText + \"' but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"...navigate the menus game knowing exactly...\" -> \"...navigate the game menus, knowing exactly...\"\nawait replaceOnce(\n  \"will be easy for any player to be able to navigate the menus game knowing exactly how to edit a specific setting or getting back to playing the game.\",\n  \"will be easy for any player to be able to navigate the game menus, knowing exactly how to edit a specific setting or getting back to playing the game.\"\n);\n\n// 2) \"...leads to a menu selector...\" -> \"...leads to a level selector...\"\nawait replaceOnce(\n  \" leads to a menu selector where you select which level to play \",\n  \" leads to a level selector where you select which level to play \"\n);\n\n// 3) \"The Options menu will have\" -> \"The options menu will have\"\nawait replaceOnce(\n  \" The Options menu will have \",\n  \" The options menu will have \"\n);\n\n// 4) Append new clause about the pause menu before the trailing period.\nawait replaceOnce(\n  \" like the video and audio settings.\",\n  \" like the video and audio settings, the menu will also be accessible through the pause menu.\"\n);\n\n// 5) \"...kind of game will feel like.\" -> \"...kind of game could feel like.\"\nawait replaceOnce(\n  \"This game will have a large influence on what our project will be like as it is a successful example of what this kind of game will feel like. \",\n  \"This game will have a large influence on what our project will be like as it is a successful example of what this kind of game could feel like. \"\n);\n", "ps1": "# Applies the textual edits described by the commit diff:\n#  1. \"navigate the menus game\" -> \"navigate the game menus,\" (reworded + comma)\n#  2. \"menu selector\" -> \"level selector\"\n#  3. \"The Options menu\" -> \"The options menu\" (capitalization)\n#  4. add \", the menu will also be accessible through the pause menu\" before the final period\n#  5. \"will feel like\" -> \"could feel like\"\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute(\n        $findText,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        $wdFindContinue, # Wrap\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        $wdReplaceAll # Replace\n    )\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n\n# 1) \"...navigate the menus game knowing exactly...\" -> \"...navigate the game menus, knowing exactly...\"\nReplace-Text `\n    \"will be easy for any player to be able to navigate the menus game knowing exactly how to edit a specific setting or getting back to playing the game.\" `\n    \"will be easy for any player to be able to navigate the game menus, knowing exactly how to edit a specific setting or getting back to playing the game.\"\n\n# 2) \"...leads to a menu selector...\" -> \"...leads to a level selector...\"\nReplace-Text `\n    \" leads to a menu selector where you select which level to play \" `\n    \" leads to a level selector where you select which level to play \"\n\n# 3) \"The Options menu will have\" -> \"The options menu will have\"\nReplace-Text `\n    \" The Options menu will have \" `\n    \" The options menu will have \"\n\n# 4) Append new clause about the pause menu before the trailing period.\nReplace-Text `\n    \" like the video and audio settings.\" `\n    \" like the video and audio settings, the menu will also be accessible through the pause menu.\"\n\n# 5) \"...kind of game will feel like.\" -> \"...kind of game could feel like.\"\nReplace-Text `\n    \"This game will have a large influence on what our project will be like as it is a successful example of what this kind of game will feel like. \" `\n    \"This game will have a large influence on what our project will be like as it is a successful example of what this kind of game could feel like. \"\n"}
